$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Date","Retail release","Central","Western District","Sheung Wan","Wan Chai",
    "Causeway Bay","Tin Hau","Happy Valley","Tai Hang","North Point","Shau Kei Wan",
    "Quarry Bay","Chai Wan","Island South","Aberdeen","Mongkok","Tsim Sha Tsui",
    "Jordan","Yau Ma Tei","Tai Kok Tsui","Tsim Sha Tsui East","Sham Shui Po",
    "Cheung Sha Wan","Mei Foo","Kowloon City","To KWa Wan","Hung Hom","Kai Tak",
    "San Po Kong","Wong Tai Sin","Kwun Tong","Ngau Tau Kok","Kowloon Bay","Yau Tong",
    "Kwai Chung","Tsing Yi","Tsuen Wan","Tuen Man","Yuen Long","Tin Shui Wai",
    "Hung Shui Kiu","Sheung Shui","Fanling","Tai Po","Sha Tin","Tai Wai","Ma On Shan",
    "Tseung Kwan O","Sai Kung","Island"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Propagate the bold/bordered header style (already on A1:F1) across the
# newly-populated header cells G1:AY1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:AY1").PasteSpecial(-4122) | Out-Null

$ws.Range("O4").Select() | Out-Null
